$wb = $excel.ActiveWorkbook

# --- Invoices sheet: update INV_ID values and selection ---
$invoices = $wb.Worksheets.Item("Invoices")
$invoices.Range("A2").Value = 1101
$invoices.Range("A3").Value = 1102

# --- Line Items sheet: update INV_ID (foreign key) values and selection ---
$lineItems = $wb.Worksheets.Item("Line Items")
$lineItems.Range("B2").Value = 1101
$lineItems.Range("B3").Value = 1101
$lineItems.Range("B4").Value = 1102
$lineItems.Range("B5").Value = 1102
$lineItems.Range("C16").Select()

# Activate Invoices sheet last so it becomes the workbook's active tab,
# and select B3 on it as the final active cell.
$invoices.Activate()
$invoices.Range("B3").Select()
